# Weekly update: insert a new Camote (Zapallo) price record ahead of the
# existing history, shifting all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 62; this pushes rows 62:92 down to 63:93
# and Excel automatically extends the sheet dimension to A1:R93.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the latest market record.
$ws.Range("A62").Value2 = 7
$ws.Range("B62").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C62").Value2 = "Ñuble"
$ws.Range("D62").Value2 = 44489
$ws.Range("E62").Value2 = 16
$ws.Range("F62").Value2 = 100112045
$ws.Range("G62").Value2 = "Zapallo"
$ws.Range("H62").Value2 = "Camote"
$ws.Range("I62").Value2 = "1a (guarda)"
$ws.Range("J62").Value2 = 160
$ws.Range("K62").Value2 = 800
$ws.Range("L62").Value2 = 900
$ws.Range("M62").Value2 = 850
$ws.Range("N62").Value2 = "`$/kilo (volumen en unidades)"
$ws.Range("O62").Value2 = "Región del Maule"
$ws.Range("P62").Value2 = 850
$ws.Range("Q62").Value2 = 1
$ws.Range("R62").Value2 = "Hortaliza"
